$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P1").Value = 0.761647890981206
$ws.Range("AM2").Value = 0.73798094802378755
$ws.Range("BN2").Value = 0.84121546814616921
$ws.Range("D3").Value = 0.7854660439301222
$ws.Range("E4").Value = 0.99767016166859102
$ws.Range("AC4").Value = 0.73961968254761024
$ws.Range("J5").Value = 0.9209185062797145
$ws.Range("AI5").Value = 0.75622846018598522
$ws.Range("AP5").Value = 0.98787672571371665
$ws.Range("G6").Value = 0.99986124680012889
$ws.Range("BL6").Value = 0.90232682372060857
$ws.Range("AE7").Value = 0.83493081581645356
$ws.Range("BA7").Value = 0.67803763656493021
$ws.Range("F8").Value = 0.79951033914849612
$ws.Range("J8").Value = 0.91209338732343681
$ws.Range("AT8").Value = 0.73032050435511708
$ws.Range("AW8").Value = 0.98740947461809969
$ws.Range("AE9").Value = 0.75854062078424156
$ws.Range("AP9").Value = 0.9954368923565341
$ws.Range("AQ9").Value = 0.84438308184582833
$ws.Range("BG10").Value = 0.97148061386419116
$ws.Range("I11").Value = 0.74286022442351096
$ws.Range("BB11").Value = 0.83906692404422678
$ws.Range("G13").Value = 0.7735864350450683
$ws.Range("Z13").Value = 0.77690325155894213
$ws.Range("AD13").Value = 0.92869590575915062
$ws.Range("Q15").Value = 0.65212291684930968
$ws.Range("AB15").Value = 0.95273849595910998
$ws.Range("AW15").Value = 0.75145413689444929
$ws.Range("N16").Value = 0.70076196208851782
$ws.Range("X16").Value = 0.65668783186657231
$ws.Range("Y16").Value = 0.67964199681754478
$ws.Range("V17").Value = 0.8812988910371754
$ws.Range("AA17").Value = 0.83932829267111608
$ws.Range("AH18").Value = 0.96563542161211813
$ws.Range("BO18").Value = 0.93858496085213849
$ws.Range("E19").Value = 0.98426777711845115
$ws.Range("L20").Value = 0.88657447921166499
$ws.Range("Q20").Value = 0.94024735030654338
$ws.Range("AK20").Value = 0.69325575194619504
$ws.Range("K21").Value = 0.95129112249702552
$ws.Range("T21").Value = 0.98021757537796494
$ws.Range("AM21").Value = 0.95116383006754446
$ws.Range("U22").Value = 0.95683132950844207
$ws.Range("K24").Value = 0.79902478155355006
$ws.Range("BL24").Value = 0.78514745061494673
$ws.Range("BO25").Value = 0.92136783315988913
$ws.Range("BP25").Value = 0.94292314439312819
$ws.Range("E26").Value = 0.77489007498136153
$ws.Range("X26").Value = 0.98844893411273382
$ws.Range("AB26").Value = 0.8214433155692944
$ws.Range("Y27").Value = 0.77852899418856736
$ws.Range("AR27").Value = 0.87599776208320601
$ws.Range("BE27").Value = 0.70202665687250754
$ws.Range("C28").Value = 0.91872889951326431
$ws.Range("V28").Value = 0.8280832011947914
$ws.Range("AA28").Value = 0.83697965664719631
$ws.Range("L29").Value = 0.87936502856502019
$ws.Range("W29").Value = 0.99717083778908677
$ws.Range("AB29").Value = 0.77047020849622161
$ws.Range("AG30").Value = 0.83069710996770896
$ws.Range("BP30").Value = 0.71810250477913184
$ws.Range("B31").Value = 0.67711564685175674
$ws.Range("AQ31").Value = 0.91867442920570386
$ws.Range("F32").Value = 0.85687921783224941
$ws.Range("M32").Value = 0.8512305005329166
$ws.Range("AH32").Value = 0.75290438541887705
$ws.Range("BC32").Value = 0.8764659519260547
$ws.Range("X33").Value = 0.71110803292334412
$ws.Range("AF33").Value = 0.95447145961468038
$ws.Range("AQ33").Value = 0.94296788592153435
$ws.Range("AT33").Value = 0.86041615020814932
$ws.Range("V34").Value = 0.66253896392606937
$ws.Range("AU34").Value = 0.94342761307774659
$ws.Range("N35").Value = 0.75424122900834578
$ws.Range("Q35").Value = 0.89564565608658175
$ws.Range("O36").Value = 0.98260423454782975
$ws.Range("AL37").Value = 0.70936651684473573
$ws.Range("G38").Value = 0.85557399258470335
$ws.Range("S38").Value = 0.8045236720153317
$ws.Range("AP38").Value = 0.85946912861694691
$ws.Range("BC39").Value = 0.99359961855528933
$ws.Range("AP40").Value = 0.79890692932235974
$ws.Range("BD40").Value = 0.74658045101029469
$ws.Range("BO40").Value = 0.96070273634188319
$ws.Range("M41").Value = 0.73349330745406849
$ws.Range("O41").Value = 0.85475630241420242
$ws.Range("L43").Value = 0.83497348120098924
$ws.Range("R43").Value = 0.90659650336568387
$ws.Range("I44").Value = 0.89243116086581742
$ws.Range("W44").Value = 0.92825183052108562
$ws.Range("BP44").Value = 0.62593988310469872
$ws.Range("O45").Value = 0.75779849286198941
$ws.Range("S45").Value = 0.83567564235507663
$ws.Range("AY45").Value = 0.86284167399691225
$ws.Range("BP45").Value = 0.9016767094988567
$ws.Range("A46").Value = 0.85451353174829359
$ws.Range("AM46").Value = 0.97710769482935811
$ws.Range("BA46").Value = 0.7327467640307801
$ws.Range("AL47").Value = 0.93991336991153607
$ws.Range("AM48").Value = 0.86338588735076827
$ws.Range("AW48").Value = 0.96084408328896798
$ws.Range("AE49").Value = 0.73698475726541313
$ws.Range("BF49").Value = 0.78419982153259316
$ws.Range("S50").Value = 0.98863770270158757
$ws.Range("BJ50").Value = 0.96561040525056319
$ws.Range("BE51").Value = 0.8449152479997647
$ws.Range("AX52").Value = 0.64719126396215809
$ws.Range("AY52").Value = 0.999983191333915
$ws.Range("BA54").Value = 0.77244178718221457
$ws.Range("BC54").Value = 0.96082296689349822
$ws.Range("BD54").Value = 0.74527184316339534
$ws.Range("BI54").Value = 0.82719099110560335
$ws.Range("BG56").Value = 0.93373880342443982
$ws.Range("AZ57").Value = 0.86323758470891288
$ws.Range("BM57").Value = 0.91103069238988599
$ws.Range("AI58").Value = 0.82009695356752044
$ws.Range("AJ58").Value = 0.97411140518326811
$ws.Range("AV58").Value = 0.70374450427887281
$ws.Range("AT60").Value = 0.91610462359715228
$ws.Range("BA60").Value = 0.79209201598502443
$ws.Range("BG60").Value = 0.89948377148978964
$ws.Range("L61").Value = 0.90399597197190928
$ws.Range("AH61").Value = 0.72895339733481923
$ws.Range("AZ61").Value = 0.52860746101088796
$ws.Range("BJ61").Value = 0.7530126066026559
$ws.Range("E62").Value = 0.91462290492247567
$ws.Range("Z63").Value = 0.66549870692595747
$ws.Range("BL63").Value = 0.70913650492442881
$ws.Range("BM63").Value = 0.90145830853643028
$ws.Range("K65").Value = 0.68394093024111391
$ws.Range("AI66").Value = 0.96439652469862225
$ws.Range("BG66").Value = 0.78040552257581886
$ws.Range("A67").Value = 0.79946757167587235
$ws.Range("AH67").Value = 0.72093505963367366
$ws.Range("E68").Value = 0.73119119786431286
